# demandes.xlsx: a new request ("Marie Martin", submitted 2025-05-28 12:52:53)
# was approved, so it is appended as row 4 of the sheet, following the same
# layout as the existing rows (2 and 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 4
$ws.Cells.Item($row, 1).Value = "20250528_125253"
$ws.Cells.Item($row, 2).Value = "2025-05-28 12:52:53"
$ws.Cells.Item($row, 3).Value = "Marie Martin"
$ws.Cells.Item($row, 4).Value = "{'chantier': 'Aluminium - Table Aluminium 02', 'urgence': 'Normal', 'date_souhaitee': '2025-05-28', 'produits': {'3119066359': {'produit': 'Crémone F8 Variable  L400 380-620', 'quantite': 9, 'emplacement': 'E2'}, '3032359406': {'produit': 'SPT/16-4,3X22-GS', 'quantite': 4, 'emplacement': 'G35'}}}"
$ws.Cells.Item($row, 5).Value = ""
$ws.Cells.Item($row, 6).Value = "Approuvée"
$ws.Cells.Item($row, 7).Value = "2025-05-28 12:58:22"
$ws.Cells.Item($row, 8).Value = "elie"
$ws.Cells.Item($row, 9).Value = "Demande approuvée et stock mis à jour"
